$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C21").Value = 289
$ws.Range("D21").Value = 248
$ws.Range("E21").Value = 41
$ws.Range("F21").Value = 71.06017191977078
